$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Atividades")

# --- Sprint 11 burndown: update the daily-progress grid (rows 8-13, cols C:I) ---
# Row 8 - "Criar e valida o formulário de ministrante"
$ws.Range("D8").ClearContents()
$ws.Range("E8").Value = 1

# Row 9 - "Modificações na cor, fonte, fundo"
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1

# Row 10 - "Refatoração do código em todas as páginas."
$ws.Range("C10").ClearContents()
$ws.Range("D10").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("G10").ClearContents()
$ws.Range("H10").Value = 1
$ws.Range("I10").ClearContents()

# Row 11 - "Criar o login do participante"
$ws.Range("D11").ClearContents()
$ws.Range("G11").Value = 2

# Row 12 - "Listar atividades disponíveis para inscrição"
$ws.Range("C12").ClearContents()
$ws.Range("D12").ClearContents()
$ws.Range("H12").Value = 1

# Row 13 - "Listar eventos disponíveis para participação."
$ws.Range("D13").ClearContents()
$ws.Range("H13").ClearContents()
$ws.Range("I13").Value = 1

$wb.Application.Calculate()

# Leave the cursor on the last-edited cell, matching the author's saved view.
[void]$ws.Range("J10").Select()
